$wb = $excel.ActiveWorkbook

# Sheet: ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 2252.8333
$ws.Range("I9").Value = 2654.125
$ws.Range("K9").Value = 2654.125
$ws.Range("M9").Value = -2485.125
$ws.Range("H55").Value = 191.26315
$ws.Range("I55").Value = 88.7
$ws.Range("J55").Value = 305.22223
$ws.Range("K55").Value = 88.7
$ws.Range("L55").Value = 305.22223
$ws.Range("M55").Value = 125.3
$ws.Range("N55").Value = -733.2222300000001
$ws.Range("H112").Value = 2081.926
$ws.Range("J112").Value = 2197.4695
$ws.Range("L112").Value = 6592.4085
$ws.Range("N112").Value = -8808.408500000001
$ws.Range("H131").Value = 2406.8572
$ws.Range("I131").Value = 2406.8572
$ws.Range("K131").Value = 7220.571599999999
$ws.Range("M131").Value = -2180.571599999999
$ws.Range("H139").Value = 143700
$ws.Range("J139").Value = 175966.67
$ws.Range("L139").Value = 175966.67
$ws.Range("N139").Value = -186246.67
$ws.Range("H141").Value = 3082.5
$ws.Range("I141").Value = 2125
$ws.Range("K141").Value = 6375
$ws.Range("M141").Value = -1195

# Sheet: ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3801.5862
$ws.Range("I2").Value = 1167.0952
$ws.Range("J2").Value = 10717.125
$ws.Range("K2").Value = 1167.0952
$ws.Range("L2").Value = 10717.125
$ws.Range("M2").Value = -1054.0952
$ws.Range("N2").Value = -10943.125
$ws.Range("H32").Value = 2322575.5
$ws.Range("I32").Value = 5487.4688
$ws.Range("K32").Value = 5487.4688
$ws.Range("M32").Value = -5200.4688
$ws.Range("H45").Value = 1545
$ws.Range("I45").Value = 1493.75
$ws.Range("K45").Value = 1493.75
$ws.Range("M45").Value = -1116.75
$ws.Range("H74").Value = 4380.1924
$ws.Range("I74").Value = 5149.5
$ws.Range("J74").Value = 3482.6667
$ws.Range("K74").Value = 5149.5
$ws.Range("L74").Value = 3482.6667
$ws.Range("M74").Value = -4275.5
$ws.Range("N74").Value = -5230.6667
$ws.Range("H77").Value = 4380.1924
$ws.Range("I77").Value = 5149.5
$ws.Range("J77").Value = 3482.6667
$ws.Range("K77").Value = 25747.5
$ws.Range("L77").Value = 17413.3335
$ws.Range("M77").Value = -21379.5
$ws.Range("N77").Value = -26149.3335
$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").ClearContents()
$ws.Range("H116").Value = 3801.5862
$ws.Range("I116").Value = 1167.0952
$ws.Range("J116").Value = 10717.125
$ws.Range("K116").Value = 1167.0952
$ws.Range("L116").Value = 10717.125
$ws.Range("M116").Value = 1126.9048
$ws.Range("N116").Value = -15305.125

# Sheet: BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3801.5862
$ws.Range("I3").Value = 1167.0952
$ws.Range("J3").Value = 10717.125
$ws.Range("K3").Value = 1167.0952
$ws.Range("L3").Value = 10717.125
$ws.Range("M3").Value = -1053.0952
$ws.Range("N3").Value = -10945.125
$ws.Range("H22").Value = 2228.4
$ws.Range("I22").Value = 232.14285
$ws.Range("J22").Value = 6886.3335
$ws.Range("K22").Value = 232.14285
$ws.Range("L22").Value = 6886.3335
$ws.Range("M22").Value = -59.14285000000001
$ws.Range("N22").Value = -7232.3335
$ws.Range("H86").Value = 2972.6858
$ws.Range("I86").Value = 1471.2307
$ws.Range("J86").Value = 7310.222
$ws.Range("K86").Value = 1471.2307
$ws.Range("L86").Value = 7310.222
$ws.Range("M86").Value = -348.2307000000001
$ws.Range("N86").Value = -9556.222
$ws.Range("H89").Value = 2972.6858
$ws.Range("I89").Value = 1471.2307
$ws.Range("J89").Value = 7310.222
$ws.Range("K89").Value = 7356.1535
$ws.Range("L89").Value = 36551.11
$ws.Range("M89").Value = -1740.1535
$ws.Range("N89").Value = -47783.11
$ws.Range("H94").Value = 6203.421
$ws.Range("I94").Value = 2965.5
$ws.Range("J94").Value = 7697.846
$ws.Range("K94").Value = 2965.5
$ws.Range("L94").Value = 7697.846
$ws.Range("M94").Value = -2514.5
$ws.Range("N94").Value = -8599.846
$ws.Range("H134").Value = 897427.4
$ws.Range("I134").Value = 1024217.6
$ws.Range("K134").Value = 3072652.8
$ws.Range("M134").Value = -3070117.8

# Sheet: CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 8342
$ws.Range("I2").Value = 10225.25
$ws.Range("J2").Value = 5831
$ws.Range("K2").Value = 10225.25
$ws.Range("L2").Value = 5831
$ws.Range("M2").Value = -10112.25
$ws.Range("N2").Value = -6057
$ws.Range("H122").Value = 11180.444
$ws.Range("I122").Value = 4902.75
$ws.Range("K122").Value = 14708.25
$ws.Range("M122").Value = -12258.25
$ws.Range("H134").Value = 55564256
$ws.Range("I134").Value = 83338310
$ws.Range("K134").Value = 250014930
$ws.Range("M134").Value = -250012395
$ws.Range("H141").Value = 990000
$ws.Range("J141").Value = 990000
$ws.Range("L141").Value = 990000
$ws.Range("N141").Value = -1000360

# Sheet: CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H129").Value = 13891493
$ws.Range("I129").Value = 829.25
$ws.Range("K129").Value = 2487.75
$ws.Range("M129").Value = 2512.25

# Sheet: GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 2260
$ws.Range("I5").Value = 2260
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 2260
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = -2148
$ws.Range("N5").ClearContents()
$ws.Range("H70").Value = 8843.5
$ws.Range("I70").Value = 9836
$ws.Range("J70").Value = 7851
$ws.Range("K70").Value = 9836
$ws.Range("L70").Value = 7851
$ws.Range("M70").Value = -9566
$ws.Range("N70").Value = -8391
$ws.Range("H73").Value = 8843.5
$ws.Range("I73").Value = 9836
$ws.Range("J73").Value = 7851
$ws.Range("K73").Value = 9836
$ws.Range("L73").Value = 7851
$ws.Range("M73").Value = -8900
$ws.Range("N73").Value = -9723
$ws.Range("H92").Value = 5083.3335
$ws.Range("J92").Value = 5083.3335
$ws.Range("L92").Value = 5083.3335
$ws.Range("N92").Value = -8827.333500000001
$ws.Range("H97").Value = 1879.5333
$ws.Range("I97").Value = 1941.5714
$ws.Range("J97").Value = 1011
$ws.Range("K97").Value = 1941.5714
$ws.Range("L97").Value = 1011
$ws.Range("M97").Value = -1445.5714
$ws.Range("N97").Value = -2003
$ws.Range("H107").Value = 844.8889
$ws.Range("I107").Value = 216.33333
$ws.Range("J107").Value = 1159.1666
$ws.Range("K107").Value = 216.33333
$ws.Range("L107").Value = 1159.1666
$ws.Range("M107").Value = 1703.66667
$ws.Range("N107").Value = -4999.1666
$ws.Range("H113").Value = 7721.706
$ws.Range("I113").Value = 3934.2727
$ws.Range("K113").Value = 3934.2727
$ws.Range("M113").Value = -1764.2727
$ws.Range("H126").Value = 38470464
$ws.Range("I126").Value = 62503564
$ws.Range("J126").Value = 17499.8
$ws.Range("K126").Value = 187510692
$ws.Range("L126").Value = 52499.39999999999
$ws.Range("M126").Value = -187508222
$ws.Range("N126").Value = -57439.39999999999
$ws.Range("H132").Value = 9326.529
$ws.Range("I132").Value = 11104.75
$ws.Range("J132").Value = 5058.8
$ws.Range("K132").Value = 33314.25
$ws.Range("L132").Value = 15176.4
$ws.Range("M132").Value = -30784.25
$ws.Range("N132").Value = -20236.4

# Sheet: LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("N2").ClearContents()
$ws.Range("H22").Value = 873.5833
$ws.Range("J22").Value = 847.25
$ws.Range("L22").Value = 847.25
$ws.Range("N22").Value = -1437.25
$ws.Range("H27").Value = 873.5833
$ws.Range("J27").Value = 847.25
$ws.Range("L27").Value = 847.25
$ws.Range("N27").Value = -1061.25
$ws.Range("H46").Value = 1807.8889
$ws.Range("I46").Value = 692.2857
$ws.Range("J46").Value = 2517.818
$ws.Range("K46").Value = 692.2857
$ws.Range("L46").Value = 2517.818
$ws.Range("M46").Value = -504.2857
$ws.Range("N46").Value = -2893.818
$ws.Range("H61").Value = 7148.8096
$ws.Range("I61").Value = 5835.625
$ws.Range("K61").Value = 5835.625
$ws.Range("M61").Value = -5633.625
$ws.Range("H82").Value = 2920.0952
$ws.Range("I82").Value = 1823.7693
$ws.Range("K82").Value = 1823.7693
$ws.Range("M82").Value = -1462.7693
$ws.Range("H85").Value = 2920.0952
$ws.Range("I85").Value = 1823.7693
$ws.Range("K85").Value = 1823.7693
$ws.Range("M85").Value = -575.7692999999999
$ws.Range("H93").Value = 2341.963
$ws.Range("I93").Value = 2347.5454
$ws.Range("J93").Value = 2317.4
$ws.Range("K93").Value = 2347.5454
$ws.Range("L93").Value = 2317.4
$ws.Range("M93").Value = -1099.5454
$ws.Range("N93").Value = -4813.4
$ws.Range("H100").Value = 2912.6667
$ws.Range("I100").Value = 3579.9167
$ws.Range("K100").Value = 3579.9167
$ws.Range("M100").Value = -3038.9167
$ws.Range("H113").Value = 7148.8096
$ws.Range("I113").Value = 5835.625
$ws.Range("K113").Value = 5835.625
$ws.Range("M113").Value = -3665.625

# Sheet: WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H95").Value = 47998.5
$ws.Range("J95").Value = 47998.5
$ws.Range("L95").Value = 47998.5
$ws.Range("N95").Value = -53490.5
$ws.Range("H126").Value = 3270.75
$ws.Range("I126").Value = 1818.2142
$ws.Range("K126").Value = 5454.642599999999
$ws.Range("M126").Value = -2984.642599999999
$ws.Range("H136").Value = 15636193
$ws.Range("I136").Value = 17249808
$ws.Range("J136").Value = 37916.668
$ws.Range("K136").Value = 51749424
$ws.Range("L136").Value = 113750.004
$ws.Range("M136").Value = -51746874
$ws.Range("N136").Value = -118850.004
